$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 5.986083666666667
$ws.Range("N2").Value = 17.958251
$ws.Range("O2").Value = 0.1400451834753423
$ws.Range("P2").Value = 0.1400451834753423
$ws.Range("Q2").Value = 54.81432470159756
$ws.Range("R2").Value = 493.3289223143781
$ws.Range("S2").Value = 0.1357558364398901
$ws.Range("T2").Value = 0.1357558364398901

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.5200697664501973
$ws.Range("P3").Value = 0.5200697664501973
$ws.Range("S3").Value = 0.5041409093799774
$ws.Range("T3").Value = 0.5041409093799774

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 14.528028
$ws.Range("N4").Value = 43.584084
$ws.Range("O4").Value = 0.3398850500744605
$ws.Range("P4").Value = 0.3398850500744605
$ws.Range("Q4").Value = 133.032561589528
$ws.Range("R4").Value = 1197.293054305752
$ws.Range("S4").Value = 0.3294749460226629
$ws.Range("T4").Value = 0.3294749460226629

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 5.986083666666667
$ws.Range("N5").Value = 17.958251
$ws.Range("O5").Value = 0.1400451834753423
$ws.Range("P5").Value = 0.1400451834753423
$ws.Range("Q5").Value = 1.731915675413445
$ws.Range("R5").Value = 15.587241078721
$ws.Range("S5").Value = 0.004289347035452157
$ws.Range("T5").Value = 0.004289347035452157

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.5200697664501973
$ws.Range("P6").Value = 0.5200697664501973
$ws.Range("Q6").Value = 6.431616985830111
$ws.Range("R6").Value = 57.884552872471
$ws.Range("S6").Value = 0.0159288570702199
$ws.Range("T6").Value = 0.0159288570702199

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 14.528028
$ws.Range("N7").Value = 43.584084
$ws.Range("O7").Value = 0.3398850500744605
$ws.Range("P7").Value = 0.3398850500744605
$ws.Range("Q7").Value = 4.203302330396001
$ws.Range("R7").Value = 37.82972097356401
$ws.Range("S7").Value = 0.01041010405179757
$ws.Range("T7").Value = 0.01041010405179757
